# The deck's theme was changed from the "Integral" (Red Violet) design to the
# stock "Office Theme" design. Concretely, the design's 12-slot theme colour
# scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) switches from the
# Integral/Red Violet palette to the default Office palette; the font scheme
# and format scheme are already identical between the two themes, so only the
# colours need to move.
#
# RGB() values below (decimal, computed as R + G*256 + B*65536) map 1:1 to the
# hex srgbClr values used by the "Office Theme" colour scheme:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$officeTheme = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeTheme[$i - 1]
}
